{"js": "// The four trailing \"liability\" sign-off lines right after\n// {{END-FOR liability}} are being dropped (the empty paragraph that used\n// to sit right before them, and the page-break paragraph that follows\n// them, both stay put):\n//   Datum: {{arzt.liabilityDateLine}}\n//   Name Patient/in: {{arzt.liabilitySignerName}}\n//   \" \" (a lone space)\n//   Unterschrift: ____________________\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Datum: {{arzt.liabilityDateLine}}\",\n  \"Name Patient/in: {{arzt.liabilitySignerName}}\",\n  \" \",\n  \"Unterschrift: ____________________\",\n]);\n\n// Only remove the specific run of paragraphs that directly follows the\n// liability FOR-loop (so we don't accidentally touch a same-text\n// paragraph elsewhere in the document).\nconst items = paragraphs.items;\nlet loopEnd = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"{{END-FOR liability}}\") {\n    loopEnd = i;\n    break;\n  }\n}\n\nif (loopEnd !== -1) {\n  // Paragraph right after the loop is the blank spacer paragraph that is\n  // kept; the run of paragraphs to delete starts right after it.\n  let start = loopEnd + 1;\n  if (start < items.length && items[start].text === \"\") {\n    start += 1;\n  }\n\n  const toDelete = [];\n  let i = start;\n  while (i < items.length && targets.has(items[i].text)) {\n    toDelete.push(items[i]);\n    i++;\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# The four trailing \"liability\" sign-off lines right after\n# {{END-FOR liability}} are being dropped (the empty paragraph that used\n# to sit right before them, and the page-break paragraph that follows\n# them, both stay put):\n#   Datum: {{arzt.liabilityDateLine}}\n#   Name Patient/in: {{arzt.liabilitySignerName}}\n#   \" \" (a lone space)\n#   Unterschrift: ____________________\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Datum: {{arzt.liabilityDateLine}}\",\n    \"Name Patient/in: {{arzt.liabilitySignerName}}\",\n    \" \",\n    \"Unterschrift: ____________________\"\n)\n\n$count = $d.Paragraphs.Count\n\n# Find the {{END-FOR liability}} marker paragraph so we only touch the\n# specific run of paragraphs that directly follows it (avoids accidentally\n# matching same-text paragraphs elsewhere in the document).\n$loopEndIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($text -eq \"{{END-FOR liability}}\") {\n        $loopEndIndex = $i\n        break\n    }\n}\n\nif ($loopEndIndex -ne -1) {\n    # Paragraph right after the loop marker is the blank spacer paragraph\n    # that is kept; the run of paragraphs to delete starts right after it.\n    $start = $loopEndIndex + 1\n    if ($start -le $count -and $d.Paragraphs.Item($start).Range.Text.TrimEnd(\"`r\") -eq \"\") {\n        $start = $start + 1\n    }\n\n    $end = $start\n    while ($end -le $count -and ($targets -contains $d.Paragraphs.Item($end).Range.Text.TrimEnd(\"`r\"))) {\n        $end = $end + 1\n    }\n    $end = $end - 1\n\n    if ($end -ge $start) {\n        $rangeStart = $d.Paragraphs.Item($start).Range.Start\n        $rangeEnd = $d.Paragraphs.Item($end).Range.End\n        $deleteRange = $d.Range($rangeStart, $rangeEnd)\n        $deleteRange.Delete()\n    }\n}\n"}
